$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Shai Gilgeous-Alexander",
    "Jordan Poole",
    "Christian Braun",
    "Kyrie Irving",
    "Zach LaVine",
    "Jalen Williams",
    "Lauri Markkanen",
    "Tobias Harris",
    "John Collins",
    "Joel Embiid",
    "Jimmy Butler",
    "RJ Barrett",
    "Terry Rozier",
    "CJ McCollum",
    "Kevin Huerter",
    "Keyonte George"
)

$positions = @(
    "PG,SG",
    "PG,SG",
    "SG,SF",
    "PG,SG",
    "SG,SF",
    "SG,SF,PF,C",
    "SF,PF",
    "SF,PF",
    "PF,C",
    "C",
    "SF,PF",
    "SG,SF,PF",
    "PG",
    "PG,SG",
    "SG,SF",
    "PG,SG"
)

$teams = @(
    "Oklahoma City Thunder",
    "Washington Wizards",
    "Denver Nuggets",
    "Dallas Mavericks",
    "Chicago Bulls",
    "Oklahoma City Thunder",
    "Utah Jazz",
    "Detroit Pistons",
    "Utah Jazz",
    "Philadelphia 76ers",
    "Miami Heat",
    "Toronto Raptors",
    "Miami Heat",
    "New Orleans Pelicans",
    "Sacramento Kings",
    "Utah Jazz"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $positions[$i]
    $ws.Cells.Item($row, 3).Value = $teams[$i]
}
